$d = $word.ActiveDocument

# Remove the sentence "The consolidation reported insufficient staff to meet
# this portion of the requirements." (including the leading space that
# precedes it), leaving the surrounding text intact:
#   "...HUD agreement. The consolidation reported insufficient staff to meet
#    this portion of the requirements. At the time..."
# becomes
#   "...HUD agreement. At the time..."
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(" The consolidation reported insufficient staff to meet this portion of the requirements.", $true, $true, $false, $false, $false, $true, 1, $false, "", 2)
